$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; temporarily unprotect so the cells can be edited,
# then re-protect once the updates are in place.
$ws.Unprotect()

# Update the confidential disclaimer: "as of 2021-03-25" -> "as of 2021-03-26"
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."
# Re-fit the row height so the multi-line text doesn't leave a stray custom height behind.
$ws.Rows(7).AutoFit()

# Update the Weight / Percent Change figures for EFA, EEM and Total rows.
$ws.Range("D2").Value = 0.8452435952845795
$ws.Range("E2").Value = 0.01295953451467846

$ws.Range("D3").Value = 0.1547564047154205
$ws.Range("E3").Value = 0.02600154083204931

$ws.Range("E4").Value = 0.0149778685226305

$ws.Protect()
